$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

# --- Estadisticos 1P ---
# Row 2 (4AEM)
$ws1.Range("D2").Value = 12
$ws1.Range("E2").Value = 9
$ws1.Range("F2").Value = 18
$ws1.Range("G2").Value = 46.15
$ws1.Range("H2").Value = 6.4

# Row 4 (4ARHM)
$ws1.Range("D4").Value = 1
$ws1.Range("E4").Value = 8
$ws1.Range("F4").Value = 15
$ws1.Range("G4").Value = 62.5
$ws1.Range("H4").Value = 6.4

# Row 5 (4BEM)
$ws1.Range("D5").Value = 0
$ws1.Range("E5").Value = 22
$ws1.Range("F5").Value = 15
$ws1.Range("G5").Value = 40.54
$ws1.Range("H5").Value = 5.8

# --- Estadisticos 2P ---
# Row 2 (4AEM)
$ws2.Range("D2").Value = 19
$ws2.Range("E2").Value = 14
$ws2.Range("F2").Value = 13
$ws2.Range("G2").Value = 33.33
$ws2.Range("H2").Value = 6.8

# Row 4 (4ARHM)
$ws2.Range("D4").Value = 6
$ws2.Range("E4").Value = 12
$ws2.Range("F4").Value = 11
$ws2.Range("G4").Value = 45.83
$ws2.Range("H4").Value = 6.6

# Row 5 (4BEM)
$ws2.Range("D5").Value = 32
$ws2.Range("E5").Value = 33
$ws2.Range("F5").Value = 4
$ws2.Range("G5").Value = 10.81
$ws2.Range("H5").Value = 6.8

# --- Estadisticos Final ---
# Row 2 (4AEM)
$ws3.Range("D2").Value = 12
$ws3.Range("E2").Value = 13
$ws3.Range("F2").Value = 14
$ws3.Range("G2").Value = 35.9
$ws3.Range("H2").Value = 6.3

# Row 4 (4ARHM)
$ws3.Range("D4").Value = 1
$ws3.Range("E4").Value = 9
$ws3.Range("F4").Value = 14
$ws3.Range("G4").Value = 58.33
$ws3.Range("H4").Value = 6.3

# Row 5 (4BEM)
$ws3.Range("D5").Value = 0
$ws3.Range("E5").Value = 23
$ws3.Range("F5").Value = 14
$ws3.Range("G5").Value = 37.84
$ws3.Range("H5").Value = 5.8

$wb.Save()
